$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    '$> *',
    '$> **',
    '$> images_*',
    '$> file[0-9]',
    '$> file[01][0-9]',
    '$> *sv',
    '$> data_*_202[0-9]',
    '$> report-*-*',
    '$> log_[a-z][a-z]',
    '$> *ch]',
    '$> module[1-9]*',
    '$> backup-*[!aeiou]',
    '$> *_summary_*',
    '$> config*[0-9]*',
    '$> test_*_v[1-9]',
    '$> *_report_[A-Z][A-Z]',
    '$> data[0-3]*_error',
    '$> *_v[0-9]0-9]0-9]',
    '$> shell_*',
    '$> user*-ctrl*',
    '$> arch[abc]ive',
    '$> *-*-*-backup',
    '$> notes_20[12][0-9]',
    '$> session*',
    '$> log-202[0-9]-[01][0-9]-*',
    '$> file*[0-9]*ch]',
    '$> data*_backup_*',
    '$> [!abc]*[xyz]*',
    '$> prefix*mid*fix*suffix',
    '$> *[^aeiou]*ch]',
    '$> report_[0-9][0-9][0-9]*',
    '$> hidden*[!0-9]',
    '$> *start*end*',
    '$> even*more*complex*pattern[0-9]',
    '$> multi[0-9][!5-7]**[a-z]',
    '$> *error*[!0-9]*',
    '$> [!a-zA-Z]**-data[!8-9]*ch]',
    '$> project_[0-9][0-9]_release*_[!x-z]',
    '$> nested*pattern*with*many[0-9]parts[abc]',
    '$> *[!0-9*]*[a-z][!A-Z]**',
    '$> customer_[A-Za-z]**[0-9]',
    '$> log-*-week[1-4]-[!5-9]*report*',
    '$> [abc][!abc]**[a-c]*x[!z]',
    '$> *ch]*_*ch]*',
    '$> prefix[0-9]*suffix*[!xyz][a-c]',
    '$> backup_[!ab]*[0-9]*_*aeiou]*z',
    '$> alpha**[0-9][A-F]*',
    '$> **[*]literal**',
    '$> abc*xyz*123*[5-7]',
    '$> [!0-9][0-9]**file[!G-M]',
    '$> long*pattern*with*many*stars*[abc]*end',
    '$> *multi*level*test[1-3]*case[!x]',
    '$> start*middle*s*t[!ae]nd',
    '$> data[2-4][A-C]*file*202[0-3]',
    '$> *_*[*][!x-z]**',
    '$> [**]*file[**]',
    '$> *ch][!hp]*',
    '$> file[0-9]*_*[a-z][!m-p]',
    '$> even*more*complex*pattern[0-9]',
    '$> session[0-9][0-9]*_*[A-Z]*final'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Reset view: remove scrolled topLeftCell and set selection to full used range
$ws.Range("A1:A60").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Wildcards updated"